$wb = $excel.ActiveWorkbook

# ALC!row12: "Don't Be So Tallow" / 'Beeswax' (Item ID 5515)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 552.2222
$ws.Cells.Item(12, 9).Value = 530.1667
$ws.Cells.Item(12, 11).Value = 530.1667
$ws.Cells.Item(12, 13).Value = -360.1667

# ALC!row69: 'Steeling the Knife, Steeling the Mind' / 'Grade 1 Mind Dissolvent' (Item ID 12616)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 14559.4
$ws.Cells.Item(69, 9).Value = 13261.75
$ws.Cells.Item(69, 10).Value = 19750
$ws.Cells.Item(69, 11).Value = 39785.25
$ws.Cells.Item(69, 12).Value = 59250
$ws.Cells.Item(69, 13).Value = -38911.25
$ws.Cells.Item(69, 14).Value = -60998

# ALC!row72: 'Surgical Substitution (L)' / 'Grade 1 Mind Dissolvent' (Item ID 12616)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 14559.4
$ws.Cells.Item(72, 9).Value = 13261.75
$ws.Cells.Item(72, 10).Value = 19750
$ws.Cells.Item(72, 11).Value = 119355.75
$ws.Cells.Item(72, 12).Value = 177750
$ws.Cells.Item(72, 13).Value = -114987.75
$ws.Cells.Item(72, 14).Value = -186486

# ARM!row2: "Ain't Got No Ingots" / 'Bronze Ingot' (Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1893.738
$ws.Cells.Item(2, 9).Value = 1523
$ws.Cells.Item(2, 10).Value = 3747.4285
$ws.Cells.Item(2, 11).Value = 1523
$ws.Cells.Item(2, 12).Value = 3747.4285
$ws.Cells.Item(2, 13).Value = -1410
$ws.Cells.Item(2, 14).Value = -3973.4285

# ARM!row61: 'Dealing with the Tough Stuff' / 'Cobalt Ingot' (Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4247.3096
$ws.Cells.Item(61, 9).Value = 4074.2896
$ws.Cells.Item(61, 11).Value = 4074.2896
$ws.Cells.Item(61, 13).Value = -3862.2896

# ARM!row96: 'The Gauntlet Is Cast' / 'High Steel Gauntlets of Fending' (Item ID 18207)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 12).Value = ""
$ws.Cells.Item(96, 14).Value = 0

# ARM!row116: 'No Scope' / 'Titanbronze Ingot' (Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1893.738
$ws.Cells.Item(116, 9).Value = 1523
$ws.Cells.Item(116, 10).Value = 3747.4285
$ws.Cells.Item(116, 11).Value = 1523
$ws.Cells.Item(116, 12).Value = 3747.4285
$ws.Cells.Item(116, 13).Value = 771
$ws.Cells.Item(116, 14).Value = -8335.4285

# ARM!row122: 'Haste for High Durium' / 'High Durium Nugget' (Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3805.25
$ws.Cells.Item(122, 9).Value = 3100
$ws.Cells.Item(122, 11).Value = 9300
$ws.Cells.Item(122, 13).Value = -6850

# ARM!row132: "Don't Bore Me, Ore Me" / 'Mountain Chromite Ingot' (Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 6314.625
$ws.Cells.Item(132, 9).Value = 6200.122
$ws.Cells.Item(132, 11).Value = 18600.366
$ws.Cells.Item(132, 13).Value = -16070.366

# ARM!row136: 'Metal with Mettle' / 'Cobalt Tungsten Ingot' (Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 4247.3096
$ws.Cells.Item(136, 9).Value = 4074.2896
$ws.Cells.Item(136, 11).Value = 12222.8688
$ws.Cells.Item(136, 13).Value = -9672.8688

# BSM!row3: 'Hells Bells' / 'Bronze Ingot' (Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1893.738
$ws.Cells.Item(3, 9).Value = 1523
$ws.Cells.Item(3, 10).Value = 3747.4285
$ws.Cells.Item(3, 11).Value = 1523
$ws.Cells.Item(3, 12).Value = 3747.4285
$ws.Cells.Item(3, 13).Value = -1409
$ws.Cells.Item(3, 14).Value = -3975.4285

# BSM!row86: 'Through Thick and Thin' / 'Adamantite Nugget' (Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 66670000
$ws.Cells.Item(86, 9).Value = 100002500
$ws.Cells.Item(86, 11).Value = 100002500
$ws.Cells.Item(86, 13).Value = -100001377

# BSM!row89: 'Piercing Eyes Deserve Piercing Shafts (L)' / 'Adamantite Nugget' (Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 66670000
$ws.Cells.Item(89, 9).Value = 100002500
$ws.Cells.Item(89, 11).Value = 500012500
$ws.Cells.Item(89, 13).Value = -500006884

# CRP!row22: 'Driving Up the Wall' / 'Elm Lumber' (Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 883.6875
$ws.Cells.Item(22, 10).Value = 1060.4286
$ws.Cells.Item(22, 12).Value = 1060.4286
$ws.Cells.Item(22, 14).Value = -1760.4286

# CRP!row31: 'Wall Not Found' / 'Walnut Lumber' (Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3441.842
$ws.Cells.Item(31, 10).Value = 7124.3
$ws.Cells.Item(31, 12).Value = 7124.3
$ws.Cells.Item(31, 14).Value = -7714.3

# CRP!row34: 'Armoires of the Rich and Famous' / 'Walnut Lumber' (Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3441.842
$ws.Cells.Item(34, 10).Value = 7124.3
$ws.Cells.Item(34, 12).Value = 7124.3
$ws.Cells.Item(34, 14).Value = -7528.3

# CRP!row58: 'You Do the Heavy Lifting' / 'Mahogany Lumber' (Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3448.6667
$ws.Cells.Item(58, 9).Value = 4113.1816
$ws.Cells.Item(58, 10).Value = 2404.4285
$ws.Cells.Item(58, 11).Value = 4113.1816
$ws.Cells.Item(58, 12).Value = 2404.4285
$ws.Cells.Item(58, 13).Value = -3910.1816
$ws.Cells.Item(58, 14).Value = -2810.4285

# CRP!row136: 'Turali Quality' / 'Dark Mahogany Lumber' (Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 3448.6667
$ws.Cells.Item(136, 9).Value = 4113.1816
$ws.Cells.Item(136, 10).Value = 2404.4285
$ws.Cells.Item(136, 11).Value = 12339.5448
$ws.Cells.Item(136, 12).Value = 7213.2855
$ws.Cells.Item(136, 13).Value = -9789.5448
$ws.Cells.Item(136, 14).Value = -12313.2855

# CUL!row86: "Let's Not Get Sappy" / 'Birch Syrup' (Item ID 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 100
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).Value = ""

# CUL!row89: 'Luxury Spillover (L)' / 'Birch Syrup' (Item ID 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(89, 8).Value = 100
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).Value = ""

# CUL!row113: "Can't Eat Just One" / 'Night Vinegar' (Item ID 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 3379.647
$ws.Cells.Item(113, 10).Value = 4224.5454
$ws.Cells.Item(113, 12).Value = 12673.6362
$ws.Cells.Item(113, 14).Value = -17013.6362

# CUL!row117: 'A Good Omen' / 'Peppered Popotoes' (Item ID 27870)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 3058.182
$ws.Cells.Item(117, 10).Value = 3314
$ws.Cells.Item(117, 12).Value = 9942
$ws.Cells.Item(117, 14).Value = -16826

# CUL!row121: 'A Cookie for Your Troubles' / 'Coffee Biscuit' (Item ID 27878)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 67108.13
$ws.Cells.Item(121, 9).Value = 271.5
$ws.Cells.Item(121, 10).Value = 111665.89
$ws.Cells.Item(121, 11).Value = 814.5
$ws.Cells.Item(121, 12).Value = 334997.67
$ws.Cells.Item(121, 13).Value = 495.5
$ws.Cells.Item(121, 14).Value = -337617.67

# CUL!row122: 'Salt of the North' / 'Northern Sea Salt' (Item ID 36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 3143.476
$ws.Cells.Item(122, 9).Value = 780.3333
$ws.Cells.Item(122, 10).Value = 4915.8335
$ws.Cells.Item(122, 11).Value = 7022.9997
$ws.Cells.Item(122, 12).Value = 44242.5015
$ws.Cells.Item(122, 13).Value = -4572.9997
$ws.Cells.Item(122, 14).Value = -49142.5015

# CUL!row123: 'Topping Up the Pot' / 'Zurek' (Item ID 36037)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(123, 8).Value = 5373.9
$ws.Cells.Item(123, 9).Value = 4934.75
$ws.Cells.Item(123, 10).Value = 5666.6665
$ws.Cells.Item(123, 11).Value = 14804.25
$ws.Cells.Item(123, 12).Value = 16999.9995
$ws.Cells.Item(123, 13).Value = -12354.25
$ws.Cells.Item(123, 14).Value = -21899.9995

# CUL!row129: 'Comfort Food' / 'Yakow Moussaka' (Item ID 36054)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 1111.2
$ws.Cells.Item(129, 10).Value = 2332.6667
$ws.Cells.Item(129, 12).Value = 6998.000100000001
$ws.Cells.Item(129, 14).Value = -16998.0001

# CUL!row131: 'The Mountain Steeped' / 'Tsai tou Vounou' (Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1776.3077
$ws.Cells.Item(131, 9).Value = 760.1818
$ws.Cells.Item(131, 10).Value = 2175.5
$ws.Cells.Item(131, 11).Value = 2280.5454
$ws.Cells.Item(131, 12).Value = 6526.5
$ws.Cells.Item(131, 13).Value = 2759.4546
$ws.Cells.Item(131, 14).Value = -16606.5

# CUL!row139: 'Najoothie' / 'Wild Banana Blend' (Item ID 44102)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 3242.5386
$ws.Cells.Item(139, 9).Value = 2144.3044
$ws.Cells.Item(139, 11).Value = 6432.9132
$ws.Cells.Item(139, 13).Value = -1292.9132

# GSM!row122: 'Awarding Academic Excellence' / 'Ametrine' (Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 7182.0713
$ws.Cells.Item(122, 9).Value = 5114.4287
$ws.Cells.Item(122, 10).Value = 9249.714
$ws.Cells.Item(122, 11).Value = 15343.2861
$ws.Cells.Item(122, 12).Value = 27749.142
$ws.Cells.Item(122, 13).Value = -12893.2861
$ws.Cells.Item(122, 14).Value = -32649.142

# GSM!row126: 'Gold Rush Order' / 'Phrygian Gold Ingot' (Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 5444
$ws.Cells.Item(126, 10).Value = 8500
$ws.Cells.Item(126, 12).Value = 25500
$ws.Cells.Item(126, 14).Value = -30440

# LTW!row40: 'Best Served Toad' / 'Toad Leather' (Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 10054.174
$ws.Cells.Item(40, 9).Value = 9656.333000000001
$ws.Cells.Item(40, 10).Value = 10309.929
$ws.Cells.Item(40, 11).Value = 9656.333000000001
$ws.Cells.Item(40, 12).Value = 10309.929
$ws.Cells.Item(40, 13).Value = -9520.333000000001
$ws.Cells.Item(40, 14).Value = -10581.929

# WVR!row5: 'Hire in the Blood' / 'Hempen Halfgloves' (Item ID 3515)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 280000
$ws.Cells.Item(5, 10).Value = 280000
$ws.Cells.Item(5, 12).Value = 280000
$ws.Cells.Item(5, 14).Value = -280224

# WVR!row81: 'Where the Dragonflies, the Net Catches' / 'Crawler Silk' (Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2186.7693
$ws.Cells.Item(81, 9).Value = 2186.7693
$ws.Cells.Item(81, 11).Value = 4373.5386
$ws.Cells.Item(81, 13).Value = -3312.5386

# WVR!row84: 'To Kill a Dragon on Nameday (L)' / 'Crawler Silk' (Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 2186.7693
$ws.Cells.Item(84, 9).Value = 2186.7693
$ws.Cells.Item(84, 11).Value = 21867.693
$ws.Cells.Item(84, 13).Value = -16563.693

# WVR!row132: 'Comfy Cabins' / 'Snow Cotton Cloth' (Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3688.9412
$ws.Cells.Item(132, 9).Value = 3721.8
$ws.Cells.Item(132, 11).Value = 11165.4
$ws.Cells.Item(132, 13).Value = -8635.400000000001

# WVR!row136: 'Weaving the Envelope' / 'Sarcenet Cloth' (Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 4358.75
$ws.Cells.Item(136, 9).Value = 4949.1875
$ws.Cells.Item(136, 10).Value = 1997
$ws.Cells.Item(136, 11).Value = 14847.5625
$ws.Cells.Item(136, 12).Value = 5991
$ws.Cells.Item(136, 13).Value = -12297.5625
$ws.Cells.Item(136, 14).Value = -11091

# WVR!row141: 'Silk for Sunperch' / 'Thunderyards Silk Coat of Casting' (Item ID 42505)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = ""
$ws.Cells.Item(141, 14).Value = 0
